$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 0.1814516129032258
$ws.Range("C2").Value = 0.5766129032258065
$ws.Range("J2").Value = 0.004032258064516129
$ws.Range("P2").Value = 0.1370967741935484
$ws.Range("S2").Value = 0.1008064516129032
$ws.Range("C3").Value = 0.03378378378378379
$ws.Range("J3").Value = 0.02702702702702703
$ws.Range("P3").Value = 0.7364864864864865
$ws.Range("S3").Value = 0.2027027027027027
$ws.Range("J4").Value = 0.06451612903225806
$ws.Range("P4").Value = 0.7419354838709677
$ws.Range("S4").Value = 0.1935483870967742
$ws.Range("B6").Value = 0.07239819004524888
$ws.Range("F6").Value = 0.06334841628959276
$ws.Range("J6").Value = 0.2398190045248869
$ws.Range("O6").Value = 0.03619909502262444
$ws.Range("Q6").Value = 0.1176470588235294
$ws.Range("R6").Value = 0.05429864253393665
$ws.Range("S6").Value = 0.416289592760181
$ws.Range("B7").Value = 0.09134615384615384
$ws.Range("D7").Value = 0.01442307692307692
$ws.Range("F7").Value = 0.0576923076923077
$ws.Range("J7").Value = 0.1538461538461539
$ws.Range("O7").Value = 0.009615384615384616
$ws.Range("Q7").Value = 0.1298076923076923
$ws.Range("R7").Value = 0.0625
$ws.Range("S7").Value = 0.4807692307692308
$ws.Range("B8").Value = 0.08816120906801007
$ws.Range("D8").Value = 0.01007556675062972
$ws.Range("E8").Value = 0.002518891687657431
$ws.Range("F8").Value = 0.06045340050377834
$ws.Range("J8").Value = 0.1309823677581864
$ws.Range("O8").Value = 0.02267002518891688
$ws.Range("Q8").Value = 0.1486146095717884
$ws.Range("R8").Value = 0.07304785894206549
$ws.Range("S8").Value = 0.4634760705289673
$ws.Range("B9").Value = 0.07407407407407407
$ws.Range("D9").Value = 0.01587301587301587
$ws.Range("F9").Value = 0.05291005291005291
$ws.Range("J9").Value = 0.07936507936507936
$ws.Range("O9").Value = 0.01058201058201058
$ws.Range("Q9").Value = 0.1746031746031746
$ws.Range("R9").Value = 0.1058201058201058
$ws.Range("S9").Value = 0.4867724867724867
$ws.Range("B10").Value = 0.09833187006145742
$ws.Range("D10").Value = 0.01755926251097454
$ws.Range("F10").Value = 0.07374890254609306
$ws.Range("J10").Value = 0.1132572431957858
$ws.Range("O10").Value = 0.01492537313432836
$ws.Range("Q10").Value = 0.1826163301141352
$ws.Range("R10").Value = 0.04916593503072871
$ws.Range("S10").Value = 0.4503950834064969
$ws.Range("G11").Value = 0.1536050156739812
$ws.Range("J11").Value = 0.07210031347962383
$ws.Range("K11").Value = 0.1849529780564263
$ws.Range("L11").Value = 0.567398119122257
$ws.Range("S11").Value = 0.0219435736677116
$ws.Range("G12").Value = 0.7704918032786885
$ws.Range("J12").Value = 0.1693989071038251
$ws.Range("K12").Value = 0.00546448087431694
$ws.Range("L12").Value = 0.01639344262295082
$ws.Range("S12").Value = 0.03825136612021858
$ws.Range("G13").Value = 0.7777777777777778
$ws.Range("J13").Value = 0.1666666666666667
$ws.Range("S13").Value = 0.05555555555555555
$ws.Range("F15").Value = 0.02538071065989848
$ws.Range("H15").Value = 0.1421319796954315
$ws.Range("I15").Value = 0.1065989847715736
$ws.Range("J15").Value = 0.3350253807106599
$ws.Range("K15").Value = 0.07106598984771574
$ws.Range("M15").Value = 0.01522842639593909
$ws.Range("O15").Value = 0.05583756345177665
$ws.Range("S15").Value = 0.2487309644670051
$ws.Range("F16").Value = 0.01818181818181818
$ws.Range("H16").Value = 0.1272727272727273
$ws.Range("I16").Value = 0.09696969696969697
$ws.Range("J16").Value = 0.4484848484848485
$ws.Range("K16").Value = 0.103030303030303
$ws.Range("O16").Value = 0.05454545454545454
$ws.Range("S16").Value = 0.1515151515151515
$ws.Range("F17").Value = 0.03179190751445087
$ws.Range("H17").Value = 0.1416184971098266
$ws.Range("I17").Value = 0.09248554913294797
$ws.Range("J17").Value = 0.4132947976878613
$ws.Range("K17").Value = 0.1040462427745665
$ws.Range("M17").Value = 0.01445086705202312
$ws.Range("O17").Value = 0.06358381502890173
$ws.Range("S17").Value = 0.138728323699422
$ws.Range("F18").Value = 0.0310077519379845
$ws.Range("H18").Value = 0.2170542635658915
$ws.Range("I18").Value = 0.08527131782945736
$ws.Range("J18").Value = 0.3875968992248062
$ws.Range("K18").Value = 0.1162790697674419
$ws.Range("M18").Value = 0.02325581395348837
$ws.Range("O18").Value = 0.0310077519379845
$ws.Range("S18").Value = 0.1085271317829457
$ws.Range("F19").Value = 0.02021660649819495
$ws.Range("H19").Value = 0.2
$ws.Range("I19").Value = 0.07725631768953069
$ws.Range("J19").Value = 0.3458483754512635
$ws.Range("K19").Value = 0.1212996389891697
$ws.Range("M19").Value = 0.01732851985559567
$ws.Range("N19").Value = 0.001444043321299639
$ws.Range("O19").Value = 0.06570397111913358
$ws.Range("S19").Value = 0.1509025270758123
